# Auto update Excel log
#
# Appends new sensor-reading rows to the PIR, Humidity and Temperature
# sheets (13 rows each, all logged at 2026-01-28 ~14:53-14:54) and grows
# each sheet's used range accordingly:
#   PIR:          A1:F70 -> A1:F83
#   Humidity:     A1:F66 -> A1:F79
#   Temperature:  A1:F66 -> A1:F79
#
# All values in this log are plain text in the source workbook (date,
# time and percentage strings are stored as text, not real Excel dates /
# numbers). Columns that look like a date/time/percentage are therefore
# force-formatted as Text ("@") immediately before the value is written,
# so Excel doesn't silently reinterpret e.g. "2026-01-28" as a date serial
# or "88.3%" as a numeric percentage.

$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append 13 new rows starting at row 71 ----
$ws = $wb.Worksheets.Item("PIR")
$PIRData = @(
        @("2026-01-28", "14:53:24", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:27", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:32", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:37", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:42", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:47", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:52", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:53:57", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:54:02", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:54:07", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:54:13", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:54:17", "14:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-28", "14:54:22", "14:00", "Bathroom", "No Motion", "Inactive")
    )
$PIRTextCols = @(1, 2, 3)
$r = 71
foreach ($row in $PIRData) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($PIRTextCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c - 1]
    }
    $r = $r + 1
}

# ---- Humidity sheet: append 13 new rows starting at row 67 ----
$ws = $wb.Worksheets.Item("Humidity")
$HumidityData = @(
        @("2026-01-28", "14:53:24", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:53:28", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:53:32", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:53:36", "14:00", "Bathroom", "87.4%", "Active"),
        @("2026-01-28", "14:53:40", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:53:44", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:53:48", "14:00", "Bathroom", "87.4%", "Active"),
        @("2026-01-28", "14:53:52", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:54:00", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:54:04", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:54:08", "14:00", "Bathroom", "87.3%", "Active"),
        @("2026-01-28", "14:54:12", "14:00", "Bathroom", "88.3%", "Active"),
        @("2026-01-28", "14:54:16", "14:00", "Bathroom", "87.4%", "Active")
    )
$HumidityTextCols = @(1, 2, 3, 5)
$r = 67
foreach ($row in $HumidityData) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($HumidityTextCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c - 1]
    }
    $r = $r + 1
}

# ---- Temperature sheet: append 13 new rows starting at row 67 ----
$ws = $wb.Worksheets.Item("Temperature")
$TemperatureData = @(
        @("2026-01-28", "14:53:24", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:28", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:32", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:36", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:40", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:44", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:48", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:53:52", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:54:00", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:54:04", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:54:08", "14:00", "Bathroom", "22.7C", "Active"),
        @("2026-01-28", "14:54:12", "14:00", "Bathroom", "22.8C", "Active"),
        @("2026-01-28", "14:54:16", "14:00", "Bathroom", "22.8C", "Active")
    )
$TemperatureTextCols = @(1, 2, 3)
$r = 67
foreach ($row in $TemperatureData) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($TemperatureTextCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c - 1]
    }
    $r = $r + 1
}

